$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166069507598877
$ws.Range("B1").Value = 2.427430868148804
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.372300386428833
$ws.Range("E1").Value = 1.234972476959229
